$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 53: hours corrected from 1 to 4
$ws.Range("C53").Value = 4

# Row 54: fill in a new journal entry (date, activity, hours)
$ws.Range("A54").Value = 43240
$ws.Range("B54").Value = "Relecture et correction rapport"
$ws.Range("C54").Value = 1.5

# Update selection to match the saved view state
$ws.Range("B55").Select()
